$wb = $excel.ActiveWorkbook

# ---- Sheet 1 ----
$ws = $wb.Worksheets.Item(1)

# Set cell values to reflect the reordered / updated rows
$ws.Range("A1").Value = "File Name"
$ws.Range("B1").Value = "zh-cn"
$ws.Range("C1").Value = "de-de"
$ws.Range("A2").Value = "ffff197dc340-9b2c-45cd-bfac-4934d7c62f95.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("A3").Value = "ffffff0deef36f-460c-460b-bdf4-32e07bd76838.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("A4").Value = "e9bcebca-39a9-41fc-b01c-754248dcd311.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# Rebuild hyperlinks: remove existing, then re-add in the same ref/target-URL order
# but with display text reflecting the new cell content (target URLs are unchanged,
# matching the fact that xl/worksheets/_rels/sheetN.xml.rels is untouched by the edit).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/e2e/e9bcebca-39a9-41fc-b01c-754248dcd311.md", "", "", "ffff197dc340-9b2c-45cd-bfac-4934d7c62f95.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/e2e/ffff197dc340-9b2c-45cd-bfac-4934d7c62f95.md", "", "", "ffffff0deef36f-460c-460b-bdf4-32e07bd76838.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/e2e/ffffff0deef36f-460c-460b-bdf4-32e07bd76838.md", "", "", "e9bcebca-39a9-41fc-b01c-754248dcd311.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/.localization-config", "", "", ".localization-config") | Out-Null

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)

# Set cell values to reflect the reordered / updated rows
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Latest Handoff File"
$ws.Range("D1").Value = "Latest Handoff Datetime"
$ws.Range("E1").Value = "Latest Target File"
$ws.Range("F1").Value = "Latest Handback File"
$ws.Range("G1").Value = "Latest Handback DateTime"
$ws.Range("H1").Value = "Handoff Reason"
$ws.Range("I1").Value = "Dependency From"
$ws.Range("A2").Value = "ffff197dc340-9b2c-45cd-bfac-4934d7c62f95.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-02 06:58:15"
$ws.Range("E2").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.md"
$ws.Range("F2").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf"
$ws.Range("G2").Value = "2016-03-02 06:59:08"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff0deef36f-460c-460b-bdf4-32e07bd76838.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-02 06:58:15"
$ws.Range("E3").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.md"
$ws.Range("F3").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf"
$ws.Range("G3").Value = "2016-03-02 06:59:08"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "e9bcebca-39a9-41fc-b01c-754248dcd311.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-02 07:04:44"
$ws.Range("E4").Value = "e9bcebca-39a9-41fc-b01c-754248dcd311.md"
$ws.Range("F4").Value = "e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.zh-cn.xlf"
$ws.Range("G4").Value = "2016-03-02 07:03:45"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Rebuild hyperlinks: remove existing, then re-add in the same ref/target-URL order
# but with display text reflecting the new cell content (target URLs are unchanged,
# matching the fact that xl/worksheets/_rels/sheetN.xml.rels is untouched by the edit).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/e2e/e9bcebca-39a9-41fc-b01c-754248dcd311.md", "", "", "ffff197dc340-9b2c-45cd-bfac-4934d7c62f95.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f0cdac184f46f90089f302f4117e618fb84dd0b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.zh-cn.xlf", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a604e18e7aec446fb86c574241fa47387ec22ca8/e2e/e9bcebca-39a9-41fc-b01c-754248dcd311.md", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d7a797ea586dc0d2a830ae0e4a58b8619633ea29/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.zh-cn.xlf", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/e2e/ffff197dc340-9b2c-45cd-bfac-4934d7c62f95.md", "", "", "ffffff0deef36f-460c-460b-bdf4-32e07bd76838.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/45d08565716534245e65888258d3629f1f858eb6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/155fe66f4466a0529e5a00c53ba0d7a67dbfbcc7/e2e/32167929-ddc0-4105-ba04-f277f96c5c5f.md", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/867b8f8b73faf0c993aea45b56f5dd28b036ac08/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/e2e/ffffff0deef36f-460c-460b-bdf4-32e07bd76838.md", "", "", "e9bcebca-39a9-41fc-b01c-754248dcd311.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/45d08565716534245e65888258d3629f1f858eb6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf", "", "", "e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/155fe66f4466a0529e5a00c53ba0d7a67dbfbcc7/e2e/32167929-ddc0-4105-ba04-f277f96c5c5f.md", "", "", "e9bcebca-39a9-41fc-b01c-754248dcd311.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/867b8f8b73faf0c993aea45b56f5dd28b036ac08/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.zh-cn.xlf", "", "", "e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/.localization-config", "", "", ".localization-config") | Out-Null

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)

# Set cell values to reflect the reordered / updated rows
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Latest Handoff File"
$ws.Range("D1").Value = "Latest Handoff Datetime"
$ws.Range("E1").Value = "Latest Target File"
$ws.Range("F1").Value = "Latest Handback File"
$ws.Range("G1").Value = "Latest Handback DateTime"
$ws.Range("H1").Value = "Handoff Reason"
$ws.Range("I1").Value = "Dependency From"
$ws.Range("A2").Value = "ffff197dc340-9b2c-45cd-bfac-4934d7c62f95.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf"
$ws.Range("D2").Value = "2016-03-02 06:58:29"
$ws.Range("E2").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.md"
$ws.Range("F2").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf"
$ws.Range("G2").Value = "2016-03-02 06:59:27"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff0deef36f-460c-460b-bdf4-32e07bd76838.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf"
$ws.Range("D3").Value = "2016-03-02 06:58:29"
$ws.Range("E3").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.md"
$ws.Range("F3").Value = "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf"
$ws.Range("G3").Value = "2016-03-02 06:59:27"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "e9bcebca-39a9-41fc-b01c-754248dcd311.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.de-de.xlf"
$ws.Range("D4").Value = "2016-03-02 07:04:55"
$ws.Range("E4").Value = "e9bcebca-39a9-41fc-b01c-754248dcd311.md"
$ws.Range("F4").Value = "e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.de-de.xlf"
$ws.Range("G4").Value = "2016-03-02 07:04:04"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Rebuild hyperlinks: remove existing, then re-add in the same ref/target-URL order
# but with display text reflecting the new cell content (target URLs are unchanged,
# matching the fact that xl/worksheets/_rels/sheetN.xml.rels is untouched by the edit).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/e2e/e9bcebca-39a9-41fc-b01c-754248dcd311.md", "", "", "ffff197dc340-9b2c-45cd-bfac-4934d7c62f95.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/21b217671d03a265b7e3d002afe8a78b97c2adf7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.de-de.xlf", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/de17fca7c5c0929d826f6302926cf17f03e6eb1c/e2e/e9bcebca-39a9-41fc-b01c-754248dcd311.md", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/29980905046c2db469be3b67a74e64c66332731f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.de-de.xlf", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/e2e/ffff197dc340-9b2c-45cd-bfac-4934d7c62f95.md", "", "", "ffffff0deef36f-460c-460b-bdf4-32e07bd76838.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a4078b19c088a4938620f108556d96064cbe2ae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/02b1f131c1e4c99f4f6299f37e06dccb193c0c95/e2e/32167929-ddc0-4105-ba04-f277f96c5c5f.md", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/52982944d5fc067fb9c3ad59a03fcf5cc5509a0e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf", "", "", "32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/e2e/ffffff0deef36f-460c-460b-bdf4-32e07bd76838.md", "", "", "e9bcebca-39a9-41fc-b01c-754248dcd311.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a4078b19c088a4938620f108556d96064cbe2ae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf", "", "", "e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/02b1f131c1e4c99f4f6299f37e06dccb193c0c95/e2e/32167929-ddc0-4105-ba04-f277f96c5c5f.md", "", "", "e9bcebca-39a9-41fc-b01c-754248dcd311.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/52982944d5fc067fb9c3ad59a03fcf5cc5509a0e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/32167929-ddc0-4105-ba04-f277f96c5c5f.10047db1deca0f979f288261c3aea2928070f3e6.de-de.xlf", "", "", "e9bcebca-39a9-41fc-b01c-754248dcd311.f9f7364cdb5dcaed988fe7c835107b35c3712b95.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/40f4742b4e92f6bcdb1c25a38bb5f6a8226d0421/.localization-config", "", "", ".localization-config") | Out-Null

